# Auto-generated edit script: applies cached market-price data refresh
# to the Leve profit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 554.6111
$ws.Range("I19").Value = 271.14285
$ws.Range("K19").Value = 271.14285
$ws.Range("M19").Value = -96.14285000000001

$ws.Range("H129").Value = 2108.611
$ws.Range("I129").Value = 1827.8334
$ws.Range("J129").Value = 2249
$ws.Range("K129").Value = 5483.5002
$ws.Range("L129").Value = 6747
$ws.Range("M129").Value = -483.5002000000004
$ws.Range("N129").Value = -16747

$ws.Range("H140").Value = 54992.5
$ws.Range("J140").Value = 54992.5
$ws.Range("L140").Value = 54992.5
$ws.Range("N140").Value = -65352.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8386.826999999999
$ws.Range("I32").Value = 5106.5835
$ws.Range("K32").Value = 5106.5835
$ws.Range("M32").Value = -4819.5835

$ws.Range("H45").Value = 1050.8334
$ws.Range("I45").Value = 1050.8334
$ws.Range("K45").Value = 1050.8334
$ws.Range("M45").Value = -673.8334

$ws.Range("H61").Value = 4241.5903
$ws.Range("I61").Value = 4271.5
$ws.Range("J61").Value = 2447
$ws.Range("K61").Value = 4271.5
$ws.Range("L61").Value = 2447
$ws.Range("M61").Value = -4059.5
$ws.Range("N61").Value = -2871

$ws.Range("H74").Value = 5347.25
$ws.Range("I74").Value = 5365.5264
$ws.Range("K74").Value = 5365.5264
$ws.Range("M74").Value = -4491.5264

$ws.Range("H77").Value = 5347.25
$ws.Range("I77").Value = 5365.5264
$ws.Range("K77").Value = 26827.632
$ws.Range("M77").Value = -22459.632

$ws.Range("H132").Value = 2465.288
$ws.Range("I132").Value = 2191.5962
$ws.Range("J132").Value = 4498.4287
$ws.Range("K132").Value = 6574.7886
$ws.Range("L132").Value = 13495.2861
$ws.Range("M132").Value = -4044.7886
$ws.Range("N132").Value = -18555.2861

$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -95060

$ws.Range("H136").Value = 4241.5903
$ws.Range("I136").Value = 4271.5
$ws.Range("J136").Value = 2447
$ws.Range("K136").Value = 12814.5
$ws.Range("L136").Value = 7341
$ws.Range("M136").Value = -10264.5
$ws.Range("N136").Value = -12441

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2762.6206
$ws.Range("I134").Value = 2147
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 6441
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -3906
$ws.Range("N134").Value = -65070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8672.797
$ws.Range("I31").Value = 5573.8
$ws.Range("J31").Value = 13192.167
$ws.Range("K31").Value = 5573.8
$ws.Range("L31").Value = 13192.167
$ws.Range("M31").Value = -5278.8
$ws.Range("N31").Value = -13782.167

$ws.Range("H34").Value = 8672.797
$ws.Range("I34").Value = 5573.8
$ws.Range("J34").Value = 13192.167
$ws.Range("K34").Value = 5573.8
$ws.Range("L34").Value = 13192.167
$ws.Range("M34").Value = -5371.8
$ws.Range("N34").Value = -13596.167

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H58").Value = 20749.834
$ws.Range("I58").Value = 26249.75
$ws.Range("J58").Value = 9750
$ws.Range("K58").Value = 26249.75
$ws.Range("L58").Value = 9750
$ws.Range("M58").Value = -26046.75
$ws.Range("N58").Value = -10156

$ws.Range("H60").Value = 25892.277
$ws.Range("I60").Value = 7018.25
$ws.Range("J60").Value = 31284.857
$ws.Range("K60").Value = 7018.25
$ws.Range("L60").Value = 31284.857
$ws.Range("M60").Value = -6507.25
$ws.Range("N60").Value = -32306.857

$ws.Range("H132").Value = 1650.3948
$ws.Range("I132").Value = 1603.0571
$ws.Range("K132").Value = 4809.1713
$ws.Range("M132").Value = -2279.1713

$ws.Range("H134").Value = 1165.7073
$ws.Range("I134").Value = 1096
$ws.Range("K134").Value = 3288
$ws.Range("M134").Value = -753

$ws.Range("H136").Value = 20749.834
$ws.Range("I136").Value = 26249.75
$ws.Range("J136").Value = 9750
$ws.Range("K136").Value = 78749.25
$ws.Range("L136").Value = 29250
$ws.Range("M136").Value = -76199.25
$ws.Range("N136").Value = -34350

$ws.Range("H140").Value = 94497.60000000001
$ws.Range("J140").Value = 94497.60000000001
$ws.Range("L140").Value = 94497.60000000001
$ws.Range("N140").Value = -104857.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 13504.556
$ws.Range("I121").Value = 4144
$ws.Range("K121").Value = 12432
$ws.Range("M121").Value = -11122

$ws.Range("H131").Value = 1691.8695
$ws.Range("I131").Value = 1576.5834
$ws.Range("J131").Value = 1817.6364
$ws.Range("K131").Value = 4729.7502
$ws.Range("L131").Value = 5452.9092
$ws.Range("M131").Value = 310.2497999999996
$ws.Range("N131").Value = -15532.9092

$ws.Range("H133").Value = 9292.362999999999
$ws.Range("J133").Value = 19999
$ws.Range("L133").Value = 59997
$ws.Range("N133").Value = -70117

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 49999.332
$ws.Range("J18").Value = 49999
$ws.Range("L18").Value = 49999
$ws.Range("N18").Value = -50585

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H38").Value = 29999
$ws.Range("J38").Value = 29999
$ws.Range("L38").Value = 29999
$ws.Range("N38").Value = -30925

$ws.Range("H132").Value = 2760.842
$ws.Range("I132").Value = 2234.8823
$ws.Range("K132").Value = 6704.646900000001
$ws.Range("M132").Value = -4174.646900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7067
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590

$ws.Range("H27").Value = 7067
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214

$ws.Range("H40").Value = 3470.8262
$ws.Range("I40").Value = 2034.5625
$ws.Range("K40").Value = 2034.5625
$ws.Range("M40").Value = -1898.5625

$ws.Range("H100").Value = 16618200
$ws.Range("I100").Value = 16618200
$ws.Range("K100").Value = 16618200
$ws.Range("M100").Value = -16617659

$ws.Range("H132").Value = 8855877
$ws.Range("I132").Value = 10963467
$ws.Range("K132").Value = 32890401
$ws.Range("M132").Value = -32887871

$ws.Range("H133").Value = 49997.6
$ws.Range("J133").Value = 49997.6
$ws.Range("L133").Value = 49997.6
$ws.Range("N133").Value = -55057.6

$ws.Range("H136").Value = 1977
$ws.Range("I136").Value = 1343.25
$ws.Range("K136").Value = 4029.75
$ws.Range("M136").Value = -1479.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2496.9333
$ws.Range("I81").Value = 2731.077
$ws.Range("J81").Value = 975
$ws.Range("K81").Value = 5462.154
$ws.Range("L81").Value = 1950
$ws.Range("M81").Value = -4401.154
$ws.Range("N81").Value = -4072

$ws.Range("H84").Value = 2496.9333
$ws.Range("I84").Value = 2731.077
$ws.Range("J84").Value = 975
$ws.Range("K84").Value = 27310.77
$ws.Range("L84").Value = 9750
$ws.Range("M84").Value = -22006.77
$ws.Range("N84").Value = -20358

$ws.Range("H132").Value = 3312.647
$ws.Range("I132").Value = 3338.7036
$ws.Range("J132").Value = 3212.1428
$ws.Range("K132").Value = 10016.1108
$ws.Range("L132").Value = 9636.428400000001
$ws.Range("M132").Value = -7486.110799999999
$ws.Range("N132").Value = -14696.4284

Write-Host "Applied market price refresh to all 8 sheets."
